$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# The new quarterly filing (1401 Q4 / "1401/12") has arrived, so every
# quarter-indexed row on the sheet rolls one column to the left (oldest
# quarter drops off column D) and the freshly reported quarter lands in M.
# Shifting via Range.Copy (instead of re-typing) preserves each cell's
# style/type exactly and avoids Excel's autoconvert-to-date parsing for
# bare "yyyy-mm-dd" text like the publish-date row.

# --- Row 8: quarter period labels ---
$ws.Range("E8:M8").Copy($ws.Range("D8:L8"))
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/12"

# --- Row 9: publish dates ---
$ws.Range("E9:M9").Copy($ws.Range("D9:L9"))
$ws.Range("M9").Value = "1402-02-13 (2)"
# The Q1-1401 report's publish-date note was revised from "(7)" to "(9)"
# after the shift landed it in column I.
$ws.Range("I9").Value = "1402-02-13 (9)"

# --- Data rows: financial figures (numbers), shift + append newest quarter ---
$ws.Range("E11:M11").Copy($ws.Range("D11:L11"))
$ws.Range("M11").Value = 4126

$ws.Range("E12:M12").Copy($ws.Range("D12:L12"))
$ws.Range("M12").Value = -2530

$ws.Range("E13:M13").Copy($ws.Range("D13:L13"))
$ws.Range("M13").Value = 1595

$ws.Range("E14:M14").Copy($ws.Range("D14:L14"))
$ws.Range("M14").Value = -197

$ws.Range("E17:M17").Copy($ws.Range("D17:L17"))
$ws.Range("M17").Value = 1399

$ws.Range("E18:M18").Copy($ws.Range("D18:L18"))
$ws.Range("M18").Value = -24

$ws.Range("E19:M19").Copy($ws.Range("D19:L19"))
$ws.Range("M19").Value = 42

$ws.Range("E20:M20").Copy($ws.Range("D20:L20"))
$ws.Range("M20").Value = 1416

$ws.Range("E21:M21").Copy($ws.Range("D21:L21"))
$ws.Range("M21").Value = 56

$ws.Range("E22:M22").Copy($ws.Range("D22:L22"))
$ws.Range("M22").Value = 1473

$ws.Range("E24:M24").Copy($ws.Range("D24:L24"))
$ws.Range("M24").Value = 1473

$ws.Range("E26:M26").Copy($ws.Range("D26:L26"))
$ws.Range("M26").Value = 440
